# FIASO_Ospedali sentinella.xlsx edit script
# Adds a new "#Ospedali" column (number of hospitals reporting) as the
# second column on the "Totale ricoveri" and "Pazienti pediatrici" sheets,
# updates a handful of data points on "Di cui in terapia intensiva", and
# restores the "Totale ricoveri" tab as the active one.

$wb = $excel.ActiveWorkbook

$wsTot  = $wb.Worksheets.Item("Totale ricoveri")
$wsTi   = $wb.Worksheets.Item("Di cui in terapia intensiva")
$wsPed  = $wb.Worksheets.Item("Pazienti pediatrici")

# ---------------------------------------------------------------------
# Sheet 1: "Totale ricoveri" - insert "#Ospedali" as new column B
# ---------------------------------------------------------------------
$wsTot.Columns("B:B").Insert()
$wsTot.Range("B1:B14").Clear()

$wsTot.Range("B1").Value = "#Ospedali"

$wsTot.Range("B2").Value = 11
$wsTot.Range("B3").Value = 16
$wsTot.Range("B4").Value = 16
$wsTot.Range("B8").Value = 21
$wsTot.Range("B9").Value = 21
$wsTot.Range("B10").Value = 21
$wsTot.Range("B11").Value = 20
$wsTot.Range("B12").Value = 20
$wsTot.Range("B13").Value = 20

# ---------------------------------------------------------------------
# Sheet 2: "Di cui in terapia intensiva" - data corrections
# ---------------------------------------------------------------------
$wsTi.Range("G12").Value = "-"
$wsTi.Range("H12").Value = 149

$wsTi.Range("E13").Value = "-"
$wsTi.Range("G13").Value = "-"
$wsTi.Range("H13").Value = "-"

$wsTi.Range("B14").Value = 195

# ---------------------------------------------------------------------
# Sheet 3: "Pazienti pediatrici" - insert "#Ospedali" as new column B
# ---------------------------------------------------------------------
$wsPed.Columns("B:B").Insert()
$wsPed.Range("B1:B12").Clear()

$wsPed.Range("B1").Value = "#Ospedali"

$wsPed.Range("B2").Value = 4
$wsPed.Range("B6").Value = 4
$wsPed.Range("B7").Value = 4
$wsPed.Range("B8").Value = 4
$wsPed.Range("B9").Value = 4
$wsPed.Range("B10").Value = 4
$wsPed.Range("B11").Value = 4
$wsPed.Range("B12").Value = 4

# ---------------------------------------------------------------------
# Restore view state: "Totale ricoveri" becomes the active/selected tab
# ---------------------------------------------------------------------
$wsPed.Activate()
$wsPed.Range("B6").Select()

$wsTi.Activate()
$wsTi.Range("B15").Select()

$wsTot.Activate()
$wsTot.Range("C14").Select()
